$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8), cloning row 2's formatting (date style on
# column A, default style elsewhere) then overwriting with the new values.
$ws.Range("A2:M2").Copy($ws.Range("A8:M8"))

$ws.Range("A8").Value = 42604.890439814815
$ws.Range("B8").Value = "Random"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 31
$ws.Range("I8").Value = 69
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 50
$ws.Range("M8").Value = 50
